# Template3.0ForMagICUpload.xlsx edit
# - sites sheet: new "dir_tilt_correction" column inserted before dir_n_samples
# - measurements sheet: "result_type"(standard)/"sequence" columns moved after
#   "citations" and demoted from required(bold) to optional(not bold); the
#   u=unknown/s=standard note text tweaked; method_codes hyperlink follows
#   its cell
# - sites sheet becomes the active tab (was locations)

$wb = $excel.ActiveWorkbook

$locations   = $wb.Worksheets.Item(1)
$sites       = $wb.Worksheets.Item(2)
$measurements = $wb.Worksheets.Item(5)

# ---------------------------------------------------------------------------
# measurements: reorder D:H from [sequence, standard, quality, method_codes,
# citations] to [quality, method_codes, citations, sequence, standard]
# ---------------------------------------------------------------------------

# drop the old method_codes hyperlink anchored on G3, capture values first
$measurements.Range("G3").Hyperlinks.Delete()

$seqNote = $measurements.Cells.Item(2, 4).Value()
$seqHdr  = $measurements.Cells.Item(3, 4).Value()

$stdHdr  = $measurements.Cells.Item(3, 5).Value()

$qualNote = $measurements.Cells.Item(2, 6).Value()
$qualHdr  = $measurements.Cells.Item(3, 6).Value()

$mcHdr  = $measurements.Cells.Item(3, 7).Value()
$citHdr = $measurements.Cells.Item(3, 8).Value()

# new column D = quality
$measurements.Cells.Item(2, 4).Value = $qualNote
$measurements.Cells.Item(3, 4).Value = $qualHdr

# new column E = method_codes
$measurements.Cells.Item(2, 5).Value = ""
$measurements.Cells.Item(3, 5).Value = $mcHdr

# new column F = citations
$measurements.Cells.Item(2, 6).Value = ""
$measurements.Cells.Item(3, 6).Value = $citHdr

# new column G = sequence (now optional/not bold)
$measurements.Cells.Item(2, 7).Value = $seqNote
$measurements.Cells.Item(3, 7).Value = $seqHdr
$measurements.Cells.Item(3, 7).Font.Bold = $false

# new column H = standard (now optional/not bold), note text updated
# (this is the first *new* shared-string text written -> lands right after
# the pre-existing strings in the table)
$measurements.Cells.Item(2, 8).Value = "u= unknown or s=standard (calibration) measurement"
$measurements.Cells.Item(3, 8).Value = $stdHdr
$measurements.Cells.Item(3, 8).Font.Bold = $false

# method_codes hyperlink now lives on E3
$measurements.Hyperlinks.Add($measurements.Range("E3"), "https://www2.earthref.org/MagIC/method-codes")

# custom width that used to sit on column G (method_codes) now sits on E
$measurements.Columns.Item(7).ColumnWidth = $measurements.Columns.Item(9).ColumnWidth
$measurements.Columns.Item(5).ColumnWidth = 15.67

$measurements.Range("H3").Select()

# ---------------------------------------------------------------------------
# sites: insert a new column R = dir_tilt_correction (with its helper note)
# ---------------------------------------------------------------------------
$sites.Columns.Item(18).Insert()

$sites.Cells.Item(2, 18).Value = "dir_tilt_correction"
$sites.Cells.Item(1, 18).Value = "Use 0 for geographic and 100 for straigraphic correction. See online data model (earthref.org/MagIC/data-models/3.0) for more options"
$sites.Columns.Item(18).ColumnWidth = 14.83

# ---------------------------------------------------------------------------
# sites becomes the active / selected sheet (was locations)
# ---------------------------------------------------------------------------
$sites.Activate()
$sites.Range("O16").Select()
